$wb = $excel.ActiveWorkbook

# --- 1. Update status text: "Ready for handoff" -> "In Translation" ---
# This shared string is used for the per-locale Status columns on every sheet.
$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = $newStatus

# --- 2. Re-fit the Status columns now that the text is shorter ---
# (target width ~13.41 chars; ColumnWidth snaps to the nearest 1/6-character
# pixel grid, so 12.5 is the closest settable value to the new narrower width)
$newWidth = 12.5
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
